$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows below had their content removed (duplicate / superseded mapping
# rows cleared out), leaving the row present but with empty cells across
# the used column range (A:O).
$rowsToClear = @(4, 5, 16, 17, 60, 77, 78, 79, 80, 81, 82, 83, 84, 85, 86, 87, 111, 112, 115, 117, 119)

foreach ($r in $rowsToClear) {
    $ws.Range("A" + $r + ":O" + $r).ClearContents()
}
